# Update cryptos list (prices / 1h volume, and swap WrappedBTC <-> BitcoinCash rows 17/18)
# as captured by the GitHub Actions refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: several "Price" (column D) values are formatted as plain-looking numbers
# (e.g. "0.3200", "1.001", "4.120") that must stay exact text, including trailing
# zeros. Forcing NumberFormat to "@" (Text) before assigning the value keeps Excel
# from re-interpreting/rounding them as numeric, then resetting the style back to
# "Normal" avoids leaving a stray explicit cell style behind.

$ws.Range('D2').Value = '30.249.35'
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').Value = '1.926.77'
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.54'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.09%  '
$ws.Range('E6').Value = '  -0.79%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3200'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.81%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '27.42'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.77%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07063'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07961'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.21%  '
$ws.Range('D13').Value = '1.931.83'
$ws.Range('E13').Value = '  +0.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.379'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.34%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.81'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.16%  '
$ws.Range('E16').Value = '  +0.82%  '
$ws.Range('B17').Value = 'BitcoinCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '258.38'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.20%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '30.263.01'
$ws.Range('E18').Value = '  -0.14%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008030'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.753'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Value = '2.184.13'
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('E22').Value = '  +0.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.846'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.528'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '166.25'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.31%  '
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('E28').Value = '  -5.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1256'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.86%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.356'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.24%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.528'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.69%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.391'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.120'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.54%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05135'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.54%  '
$ws.Range('E35').Value = '  +4.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7447'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.77%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.765'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.74%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01959'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.36%  '
$ws.Range('E39').Value = '  -1.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.358'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.53%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4496'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.986'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8455'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.11%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.001'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '100.66'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.722'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.421'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.79%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '36.55'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06104'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.50%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4192'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.94%  '
